{"js": "// 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nconst dateResults = context.document.body.search(\"September 19, 2025\", { matchCase: true });\ndateResults.load(\"text\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n    dateResults.items[0].insertText(\"September 21, 2025\", \"Replace\");\n}\n\n// 2) Split the single-line mailing address paragraph\n//    \"919 Story Road, San Jose CA 95122\" into two paragraphs:\n//    \"919 Story Road\" and a new paragraph \"San Jose, CA 95122\".\n//    (Only the standalone body paragraph is touched; the identical text\n//    inside the \"PROPERTY ADDRESS\" table cell must stay untouched.)\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nlet addressParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text === \"919 Story Road, San Jose CA 95122\") {\n        addressParagraph = paragraphs.items[i];\n        break;\n    }\n}\nif (addressParagraph) {\n    addressParagraph.insertText(\"919 Story Road\", \"Replace\");\n    addressParagraph.insertParagraph(\"San Jose, CA 95122\", \"After\");\n    await context.sync();\n}\n\n// 3) Remove the empty \"No Spacing\" paragraph that immediately follows the\n//    \"...Board of Directors\" paragraph\nconst paragraphs2 = context.document.body.paragraphs;\nparagraphs2.load(\"text\");\nawait context.sync();\n\nlet boardIndex = -1;\nfor (let i = 0; i < paragraphs2.items.length; i++) {\n    if (paragraphs2.items[i].text.endsWith(\"Board of Directors\")) {\n        boardIndex = i;\n        break;\n    }\n}\nif (boardIndex !== -1 && boardIndex + 1 < paragraphs2.items.length) {\n    const afterBoard = paragraphs2.items[boardIndex + 1];\n    afterBoard.load(\"text\");\n    await context.sync();\n    if (afterBoard.text === \"\") {\n        afterBoard.delete();\n        await context.sync();\n    }\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Update the letter date: \"September 19, 2025\" -> \"September 21, 2025\"\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq \"September 19, 2025\") {\n        $p.Range.Text = \"September 21, 2025\"\n        break\n    }\n}\n\n# 2) Split the single-line mailing address paragraph\n#    \"919 Story Road, San Jose CA 95122\" into two paragraphs:\n#    \"919 Story Road\" and a new paragraph \"San Jose, CA 95122\"\n$target = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -eq \"919 Story Road, San Jose CA 95122\") {\n        $target = $p\n        break\n    }\n}\nif ($target -ne $null) {\n    $target.Range.Text = \"919 Story Road\"\n    $target.Range.InsertParagraphAfter()\n    $newPara = $target.Next()\n    $newPara.Range.Text = \"San Jose, CA 95122\"\n}\n\n# 3) Remove the empty \"NoSpacing\" paragraph that immediately follows the\n#    \"...Board of Directors\" paragraph\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.TrimEnd(\"`r`a\") -match \"Board of Directors$\") {\n        $after = $p.Next()\n        if ($after -ne $null -and $after.Range.Text.TrimEnd(\"`r`a\") -eq \"\") {\n            $after.Range.Delete()\n        }\n        break\n    }\n}\n"}
